# Anonymize mock transaction data ("changed mock data and tests to be less
# private"): the free-text description (column C) and payee/business name
# (column H) columns contained real vendor/person/address data. Replace them
# with generic "<Category> Payment for Business N" / "Business N" placeholders.
# Columns A, B, D, E, F, G, I, J (dates, amounts, formulas, categories) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Utilities Payment for Business 4"
$ws.Range("H1").Value = "Business 4"

$ws.Range("C2").Value = "Utiltiies Payment for Business 1"
$ws.Range("H2").Value = "Business 4"

$ws.Range("C3").Value = "Taxes Payment for Business 1"
$ws.Range("H3").Value = "Business 1"

$ws.Range("C4").Value = "Repairs Payment for Business 2"
$ws.Range("H4").Value = "Business 2"

$ws.Range("C5").Value = "Cleaning & Maintenance Payment for Business 2"
$ws.Range("H5").Value = "Business 2"

$ws.Range("C6").Value = "Meal Payment for Business 4"
$ws.Range("H6").Value = "Business 4"

$ws.Range("C7").Value = "Other Payment for Business 3"
$ws.Range("H7").Value = "Business 3"
